$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old leading index column (A) -- everything shifts one column
#    to the left (old B -> A, old C -> B, ... old J -> I).
# ---------------------------------------------------------------------------
$ws.Columns("A").Delete()

# ---------------------------------------------------------------------------
# 2. Replace the boolean "which constraint is relaxed" flags (now in C:E)
#    with textual yes/no values, and make the "#Found" column (F, was G)
#    a constant 10 for every scenario re-run under the new LaTeX numbers.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "no"
$ws.Range("D2").Value = "no"
$ws.Range("E2").Value = "no"

$ws.Range("C3").Value = "no"
$ws.Range("D3").Value = "no"
$ws.Range("E3").Value = "yes"

$ws.Range("C4").Value = "no"
$ws.Range("D4").Value = "yes"
$ws.Range("E4").Value = "no"

$ws.Range("C5").Value = "yes"
$ws.Range("D5").Value = "no"
$ws.Range("E5").Value = "no"

$ws.Range("F2:F5").Value = 10

# ---------------------------------------------------------------------------
# 3. Updated Gap / Runtime / Time_Prepro numbers from the re-run experiment.
# ---------------------------------------------------------------------------
$ws.Range("H2").Value = 3601.4
$ws.Range("I2").Value = 0

$ws.Range("H3").Value = 2863.76
$ws.Range("I3").Value = 6.62

$ws.Range("H4").Value = 2920.99
$ws.Range("I4").Value = 5.31

$ws.Range("G5").Value = 0.29
$ws.Range("H5").Value = 2643.25
$ws.Range("I5").Value = 7.35

# ---------------------------------------------------------------------------
# 4. A brand-new scenario (all three flags "yes") reported as row 6.
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = 18
$ws.Range("C6").Value = "yes"
$ws.Range("D6").Value = "yes"
$ws.Range("E6").Value = "yes"
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 0.03
$ws.Range("H6").Value = 1453.59
$ws.Range("I6").Value = 7.55

Write-Host "data written"

# ---------------------------------------------------------------------------
# 5. Merge the n_N / n_B instance-size columns down the five data rows.
# ---------------------------------------------------------------------------
$ws.Range("A2:A6").Merge()
$ws.Range("B2:B6").Merge()

# ---------------------------------------------------------------------------
# 6. Re-box the table: thin border all around + a vertical divider before
#    the "#Found" column, which separates the instance/config columns from
#    the solver-result columns. Centre everything horizontally, and centre
#    vertically the two merged instance-size columns + the results on the
#    last row.
# ---------------------------------------------------------------------------
$all = $ws.Range("A1:I6")
$all.HorizontalAlignment = -4108
$all.Borders.LineStyle = -4142   # clear the old leftover full-box borders first

$ws.Range("A1:I1").Borders.Item(9).LineStyle = 1
$ws.Range("A2:I2").Borders.Item(8).LineStyle = 1
$ws.Range("A6:I6").Borders.Item(9).LineStyle = 1
$ws.Range("F1:F6").Borders.Item(7).LineStyle = 1

$ws.Range("A2:B6").VerticalAlignment = -4108
$ws.Range("H6:I6").VerticalAlignment = -4108

Write-Host "styled"

# ---------------------------------------------------------------------------
# 7. Column widths move with the data (old D/E/F/J -> new C/D/E/I).
# ---------------------------------------------------------------------------
$ws.Range("I19").Select()

